$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins/Losses/Ties), matching the formatting used by
# the rest of the header row (bold, centered, top-aligned, thin border box)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fill in the team record (Wins/Losses/Ties) for every data row
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 96   # AD
    $ws.Cells.Item($r, 31).Value = 66   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
